# Add two new columns "I0" (I) and "IF" (J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they keep the same bold/border/centered style (s="1").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-29) ---
$data = @(
    @(8, 9),
    @(5, 6),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(6, 7),
    @(8, 9),
    @(8, 9),
    @(7, 8),
    @(7, 8),
    @(12, 12),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(4, 4),
    @(6, 6),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
